$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell content updates -------------------------------------------------
# Column J (MessageType) rows 2, 6, 7: "JSON" -> "JSONMessageType"
# Column H (Event) rows 2, 6, 7: "TEST" -> "DemoEvent"
# (J is written first so the new shared-string table ends up in the same
#  order as the target workbook: JSONMessageType then DemoEvent.)
$ws.Range("J2").Value = "JSONMessageType"
$ws.Range("J6").Value = "JSONMessageType"
$ws.Range("J7").Value = "JSONMessageType"

$ws.Range("H2").Value = "DemoEvent"
$ws.Range("H6").Value = "DemoEvent"
$ws.Range("H7").Value = "DemoEvent"

# --- Column widths ---------------------------------------------------------
# Column H (8): width 6 -> ~10.5703125
$ws.Columns.Item(8).ColumnWidth = 9.7
# Column J (10): previously default width -> 18
$ws.Columns.Item(10).ColumnWidth = 17.15

# --- View / selection -------------------------------------------------------
# Scroll the sheet so column F is the left-most visible column, then select I19
$excel.ActiveWindow.ScrollColumn = 6
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("I19").Select()
